$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$wsALC.Range("H4").Value = 122.57143
$wsALC.Range("I4").Value = 109.333336
$wsALC.Range("K4").Value = 109.333336
$wsALC.Range("M4").Value = 4.666663999999997
$wsALC.Range("H26").Value = 1577
$wsALC.Range("J26").Value = 1577
$wsALC.Range("L26").Value = 1577
$wsALC.Range("N26").Value = -2265
$wsALC.Range("H51").Value = 4650
$wsALC.Range("J51").Value = 3000.3333
$wsALC.Range("L51").Value = 3000.3333
$wsALC.Range("N51").Value = -3968.3333
$wsALC.Range("H62").Value = 8589.111000000001
$wsALC.Range("I62").Value = 9333.333000000001
$wsALC.Range("J62").Value = 8217
$wsALC.Range("K62").Value = 9333.333000000001
$wsALC.Range("L62").Value = 8217
$wsALC.Range("M62").Value = -8709.333000000001
$wsALC.Range("N62").Value = -9465
$wsALC.Range("H65").Value = 8589.111000000001
$wsALC.Range("I65").Value = 9333.333000000001
$wsALC.Range("J65").Value = 8217
$wsALC.Range("K65").Value = 46666.665
$wsALC.Range("L65").Value = 41085
$wsALC.Range("M65").Value = -43546.665
$wsALC.Range("N65").Value = -47325
$wsALC.Range("H88").Value = 1487.2858
$wsALC.Range("J88").Value = 1651.3334
$wsALC.Range("L88").Value = 1651.3334
$wsALC.Range("N88").Value = -2463.3334
$wsALC.Range("H91").Value = 1487.2858
$wsALC.Range("J91").Value = 1651.3334
$wsALC.Range("L91").Value = 1651.3334
$wsALC.Range("N91").Value = -4459.3334
$wsALC.Range("H98").Value = 647.5
$wsALC.Range("I98").Value = 624
$wsALC.Range("K98").Value = 624
$wsALC.Range("M98").Value = 874
$wsALC.Range("H122").Value = 647.5
$wsALC.Range("I122").Value = 624
$wsALC.Range("K122").Value = 1872
$wsALC.Range("M122").Value = 578
$wsALC.Range("H132").Value = 37040496
$wsALC.Range("I132").Value = 43482068
$wsALC.Range("K132").Value = 130446204
$wsALC.Range("M132").Value = -130443674
$wsALC.Range("H135").Value = 15628778
$wsALC.Range("I135").Value = 676.7778
$wsALC.Range("J135").Value = 100020530
$wsALC.Range("K135").Value = 6091.000199999999
$wsALC.Range("L135").Value = 900184770
$wsALC.Range("M135").Value = -3556.000199999999
$wsALC.Range("N135").Value = -900189840
$wsALC.Range("H138").Value = 3971.6667
$wsALC.Range("I138").Value = 2870.4167
$wsALC.Range("J138").Value = 5440
$wsALC.Range("K138").Value = 8611.250100000001
$wsALC.Range("L138").Value = 16320
$wsALC.Range("M138").Value = -3471.250100000001
$wsALC.Range("N138").Value = -26600
$wsALC.Range("H141").Value = 2134.5405
$wsALC.Range("I141").Value = 1837.5883
$wsALC.Range("K141").Value = 5512.7649
$wsALC.Range("M141").Value = -332.7649000000001

# --- ARM ---
$wsARM.Range("H61").Value = 18523878
$wsARM.Range("I61").Value = 22227540
$wsARM.Range("K61").Value = 22227540
$wsARM.Range("M61").Value = -22227328
$wsARM.Range("H136").Value = 18523878
$wsARM.Range("I136").Value = 22227540
$wsARM.Range("K136").Value = 66682620
$wsARM.Range("M136").Value = -66680070

# --- BSM ---
$wsBSM.Range("H22").Value = 1124.875
$wsBSM.Range("I22").Value = 1274.25
$wsBSM.Range("J22").Value = 975.5
$wsBSM.Range("K22").Value = 1274.25
$wsBSM.Range("L22").Value = 975.5
$wsBSM.Range("M22").Value = -1101.25
$wsBSM.Range("N22").Value = -1321.5
$wsBSM.Range("H134").Value = 5683.4165
$wsBSM.Range("I134").Value = 5929.2383
$wsBSM.Range("J134").Value = 3962.6667
$wsBSM.Range("K134").Value = 17787.7149
$wsBSM.Range("L134").Value = 11888.0001
$wsBSM.Range("M134").Value = -15252.7149
$wsBSM.Range("N134").Value = -16958.0001

# --- CRP ---
$wsCRP.Range("H7").Value = 81.166664
$wsCRP.Range("J7").Value = 42
$wsCRP.Range("L7").Value = 42
$wsCRP.Range("N7").Value = -268
$wsCRP.Range("H22").Value = 151.36842
$wsCRP.Range("I22").Value = 145.06667
$wsCRP.Range("K22").Value = 145.06667
$wsCRP.Range("M22").Value = 204.93333
$wsCRP.Range("H31").Value = 4571.1274
$wsCRP.Range("I31").Value = 2364.6316
$wsCRP.Range("K31").Value = 2364.6316
$wsCRP.Range("M31").Value = -2069.6316
$wsCRP.Range("H34").Value = 4571.1274
$wsCRP.Range("I34").Value = 2364.6316
$wsCRP.Range("K34").Value = 2364.6316
$wsCRP.Range("M34").Value = -2162.6316
$wsCRP.Range("H52").Value = 22750
$wsCRP.Range("I52").Value = 8000
$wsCRP.Range("J52").Value = 26964.285
$wsCRP.Range("K52").Value = 8000
$wsCRP.Range("L52").Value = 26964.285
$wsCRP.Range("M52").Value = -7706
$wsCRP.Range("N52").Value = -27552.285
$wsCRP.Range("H58").Value = 13284.357
$wsCRP.Range("I58").Value = 1508.1666
$wsCRP.Range("J58").Value = 22116.5
$wsCRP.Range("K58").Value = 1508.1666
$wsCRP.Range("L58").Value = 22116.5
$wsCRP.Range("M58").Value = -1305.1666
$wsCRP.Range("N58").Value = -22522.5
$wsCRP.Range("H132").Value = 27029284
$wsCRP.Range("I132").Value = 38463076
$wsCRP.Range("J132").Value = 3955.5454
$wsCRP.Range("K132").Value = 115389228
$wsCRP.Range("L132").Value = 11866.6362
$wsCRP.Range("M132").Value = -115386698
$wsCRP.Range("N132").Value = -16926.6362
$wsCRP.Range("H134").Value = 55556570
$wsCRP.Range("I134").Value = 62500980
$wsCRP.Range("J134").Value = 1257
$wsCRP.Range("K134").Value = 187502940
$wsCRP.Range("L134").Value = 3771
$wsCRP.Range("M134").Value = -187500405
$wsCRP.Range("N134").Value = -8841
$wsCRP.Range("H136").Value = 13284.357
$wsCRP.Range("I136").Value = 1508.1666
$wsCRP.Range("J136").Value = 22116.5
$wsCRP.Range("K136").Value = 4524.4998
$wsCRP.Range("L136").Value = 66349.5
$wsCRP.Range("M136").Value = -1974.4998
$wsCRP.Range("N136").Value = -71449.5

# --- CUL ---
$wsCUL.Range("H5").Value = 1288.8182
$wsCUL.Range("J5").Value = 2277.7144
$wsCUL.Range("L5").Value = 6833.1432
$wsCUL.Range("N5").Value = -7057.1432
$wsCUL.Range("H63").Value = 3122.8
$wsCUL.Range("I63").Value = 1671.4286
$wsCUL.Range("K63").Value = 5014.2858
$wsCUL.Range("M63").Value = -4265.2858
$wsCUL.Range("H66").Value = 3122.8
$wsCUL.Range("I66").Value = 1671.4286
$wsCUL.Range("K66").Value = 15042.8574
$wsCUL.Range("M66").Value = -11298.8574
$wsCUL.Range("H70").Value = 3703.0908
$wsCUL.Range("I70").Value = 1412.5
$wsCUL.Range("K70").Value = 4237.5
$wsCUL.Range("M70").Value = -3922.5
$wsCUL.Range("H73").Value = 3703.0908
$wsCUL.Range("I73").Value = 1412.5
$wsCUL.Range("K73").Value = 4237.5
$wsCUL.Range("M73").Value = -3145.5
$wsCUL.Range("H75").Value = 475
$wsCUL.Range("I75").Value = 0
$wsCUL.Range("J75").Value = 475
$wsCUL.Range("K75").Value = 0
$wsCUL.Range("L75").Value = 1425
$wsCUL.Range("M75").ClearContents()
$wsCUL.Range("N75").Value = -3421
$wsCUL.Range("H78").Value = 475
$wsCUL.Range("I78").Value = 0
$wsCUL.Range("J78").Value = 475
$wsCUL.Range("K78").Value = 0
$wsCUL.Range("L78").Value = 4275
$wsCUL.Range("M78").ClearContents()
$wsCUL.Range("N78").Value = -14259
$wsCUL.Range("H87").Value = 22350.75
$wsCUL.Range("I87").Value = 9600
$wsCUL.Range("K87").Value = 28800
$wsCUL.Range("M87").Value = -27552
$wsCUL.Range("H90").Value = 22350.75
$wsCUL.Range("I90").Value = 9600
$wsCUL.Range("K90").Value = 86400
$wsCUL.Range("M90").Value = -80160
$wsCUL.Range("H103").Value = 3241.5
$wsCUL.Range("J103").Value = 8352
$wsCUL.Range("L103").Value = 25056
$wsCUL.Range("N103").Value = -26814
$wsCUL.Range("H113").Value = 806.913
$wsCUL.Range("J113").Value = 879.26666
$wsCUL.Range("L113").Value = 2637.79998
$wsCUL.Range("N113").Value = -6977.79998
$wsCUL.Range("H131").Value = 720.8099999999999
$wsCUL.Range("J131").Value = 759.9101000000001
$wsCUL.Range("L131").Value = 2279.7303
$wsCUL.Range("N131").Value = -12359.7303
$wsCUL.Range("H135").Value = 1288.8182
$wsCUL.Range("J135").Value = 2277.7144
$wsCUL.Range("L135").Value = 20499.4296
$wsCUL.Range("N135").Value = -25569.4296

# --- GSM ---
$wsGSM.Range("H122").Value = 4315.625
$wsGSM.Range("I122").Value = 4605.357
$wsGSM.Range("K122").Value = 13816.071
$wsGSM.Range("M122").Value = -11366.071
$wsGSM.Range("H126").Value = 3748
$wsGSM.Range("J126").Value = 5511.1113
$wsGSM.Range("L126").Value = 16533.3339
$wsGSM.Range("N126").Value = -21473.3339
$wsGSM.Range("H132").Value = 8506968
$wsGSM.Range("I132").Value = 18153430
$wsGSM.Range("J132").Value = 66313.875
$wsGSM.Range("K132").Value = 54460290
$wsGSM.Range("L132").Value = 198941.625
$wsGSM.Range("M132").Value = -54457760
$wsGSM.Range("N132").Value = -204001.625
$wsGSM.Range("H140").Value = 60700
$wsGSM.Range("J140").Value = 60700
$wsGSM.Range("L140").Value = 60700
$wsGSM.Range("N140").Value = -71060

# --- LTW ---
$wsLTW.Range("H9").Value = 600
$wsLTW.Range("I9").Value = 400
$wsLTW.Range("J9").Value = 800
$wsLTW.Range("K9").Value = 400
$wsLTW.Range("L9").Value = 800
$wsLTW.Range("M9").Value = -176
$wsLTW.Range("N9").Value = -1248
$wsLTW.Range("H22").Value = 4533.6665
$wsLTW.Range("I22").Value = 10001
$wsLTW.Range("K22").Value = 10001
$wsLTW.Range("M22").Value = -9706
$wsLTW.Range("H27").Value = 4533.6665
$wsLTW.Range("I27").Value = 10001
$wsLTW.Range("K27").Value = 10001
$wsLTW.Range("M27").Value = -9894
$wsLTW.Range("H40").Value = 129111
$wsLTW.Range("I40").Value = 163714.14
$wsLTW.Range("J40").Value = 8000
$wsLTW.Range("K40").Value = 163714.14
$wsLTW.Range("L40").Value = 8000
$wsLTW.Range("M40").Value = -163578.14
$wsLTW.Range("N40").Value = -8272

# --- WVR ---
$wsWVR.Range("H100").Value = 333.8889
$wsWVR.Range("I100").Value = 421
$wsWVR.Range("K100").Value = 842
$wsWVR.Range("M100").Value = -301
$wsWVR.Range("H126").Value = 1840.5
$wsWVR.Range("J126").Value = 2967.5
$wsWVR.Range("L126").Value = 8902.5
$wsWVR.Range("N126").Value = -13842.5
$wsWVR.Range("H136").Value = 3375.2812
$wsWVR.Range("I136").Value = 1750.1818
$wsWVR.Range("K136").Value = 5250.5454
$wsWVR.Range("M136").Value = -2700.5454
